# Fruta / hortaliza, semanal
# Insert a new weekly record as row 188 in the "Ciruela" price sheet,
# shifting the existing rows 188-243 down to 189-244.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 188 (pushes old rows 188-243 -> 189-244)
$ws.Rows.Item(188).Insert()

# Populate the new row 188 with the new weekly price record
$ws.Cells.Item(188, 1).Value2  = 10
$ws.Cells.Item(188, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(188, 3).Value2  = "La Araucanía"
$ws.Cells.Item(188, 4).Value2  = 44900
$ws.Cells.Item(188, 5).Value2  = 9
$ws.Cells.Item(188, 6).Value2  = "Fruta"
$ws.Cells.Item(188, 7).Value2  = 100103
$ws.Cells.Item(188, 8).Value2  = "Frutos de hueso (carozo)"
$ws.Cells.Item(188, 9).Value2  = 100103002
$ws.Cells.Item(188, 10).Value2 = "Ciruela"
$ws.Cells.Item(188, 11).Value2 = "Angeleno"
$ws.Cells.Item(188, 12).Value2 = "Primera"
$ws.Cells.Item(188, 13).Value2 = 100
$ws.Cells.Item(188, 14).Value2 = 28000
$ws.Cells.Item(188, 15).Value2 = 28000
$ws.Cells.Item(188, 16).Value2 = 28000
$ws.Cells.Item(188, 17).Value2 = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(188, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(188, 19).Value2 = 1556
$ws.Cells.Item(188, 20).Value2 = 18
